# Apply the edit described by the diff:
#  - Rename sheet "DataSet" -> "Credentials"
#  - Make the Credentials sheet (formerly DataSet) the active tab, replacing
#    ShoppingItems as the active/selected sheet
#  - Move the selection on the Credentials sheet to C2 (was D6)
#  - ShoppingItems loses tabSelected and keeps its existing selection (B17)

$wb = $excel.ActiveWorkbook

$wsData = $wb.Worksheets.Item("DataSet")
$wsData.Name = "Credentials"

# Activating the sheet and selecting C2 makes Credentials the active/selected
# tab (moving tabSelected off ShoppingItems) and updates the selection.
$wsData.Activate()
$wsData.Range("C2").Select() | Out-Null
